# Apply the changes described by the commit "Ajout de la video finale et
# modification de la lettre de motivation":
#   1. Right-align the four header/contact paragraphs (phone, email,
#      address, GUADELOUPE).
#   2. Shrink the body-letter font size from 16pt (sz/szCs 32) down to
#      15pt (sz/szCs 30) across the three main paragraphs of the letter
#      (runs + paragraph mark).

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    # Locate the paragraph containing $needle and return a Word.Paragraph.
    $r = $doc.Content
    $found = $r.Find.Execute($needle, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    return $r.Paragraphs(1)
}

# --- 1. Right-align the contact-info paragraphs -----------------------
$contactAnchors = @(
    "0690.05.06.55",
    "mathys971i22@gmail.com",
    "Rue Joseph ANTENOR, CALVAIRE 971",
    "GUADELOUPE"
)

foreach ($anchor in $contactAnchors) {
    $para = Find-ParagraphByText $d $anchor
    if ($para -ne $null) {
        $para.Alignment = 2   # wdAlignParagraphRight
    }
}

# --- 2. Shrink the letter body font size 16pt -> 15pt ------------------
$bodyAnchors = @(
    "Bonjour, je ",
    "C'est pour cela que je me tourne vers ",
    "J'espère que cette lettre débouchera vers "
)

foreach ($anchor in $bodyAnchors) {
    $para = Find-ParagraphByText $d $anchor
    if ($para -ne $null) {
        $para.Range.Font.Size = 15
        $para.Range.Font.SizeBi = 15
    }
}
